## Insert a new weekly price-report row for "Murcott" mandarinas at
## Vega Monumental Concepción, shifting the existing rows 85-168 down to
## 86-169 (dimension grows from A1:T168 to A1:T169).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 85; Excel's InsertShift pushes
# every row at/after 85 down by one and carries formatting from the row
# above, matching the workbook's existing style (date format on column D).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the reported values.
$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value = 44874
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100102
$ws.Range("H85").Value = "Cítricos"
$ws.Range("I85").Value = 100102004
$ws.Range("J85").Value = "Mandarina"
$ws.Range("K85").Value = "Murcott"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 200
$ws.Range("N85").Value = 8000
$ws.Range("O85").Value = 9000
$ws.Range("P85").Value = 8500
$ws.Range("Q85").Value = "$/bandeja 18 kilos"
$ws.Range("R85").Value = "Región de O'Higgins"
$ws.Range("S85").Value = 472
$ws.Range("T85").Value = 18
